# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") holds strikeout totals that were recalculated; update the
# values for rows 2-30 and 32 (row 31 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 3
    4  = 4
    5  = 7
    6  = 7
    7  = 6
    8  = 8
    9  = 11
    10 = 5
    11 = 4
    12 = 2
    13 = 6
    14 = 3
    15 = 9
    16 = 0
    17 = 10
    18 = 5
    19 = 4
    20 = 6
    21 = 5
    22 = 6
    23 = 5
    24 = 7
    25 = 6
    26 = 5
    27 = 3
    28 = 3
    29 = 4
    30 = 2
    32 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
